# Apply the Alvearie FHIR IG "savings-type" StructureDefinition metadata
# refresh (version bump 5.0.0 -> 6.0.0, regenerated date, publisher/
# jurisdiction populated, duplicate "Contact" metadata row removed, and the
# root Extension element's Short/Definition text updated) as captured by
# the commit diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet "Metadata" (first sheet) ----
$meta = $wb.Worksheets.Item(1)

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: refreshed publication timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was previously blank
$meta.Range("B9").Value = "Alvearie Team"

# The old row 10 ("Contact" / "No display for ContactDetail") becomes the
# Jurisdiction row.
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row;
# remove it entirely so everything below shifts up one row.
$meta.Rows.Item(11).Delete()

# ---- Sheet "Elements" (second sheet) ----
$elements = $wb.Worksheets.Item(2)

# Root Extension row: Short/Definition switch from the generic placeholder
# text to the extension-specific text.
$elements.Range("K2").Value = "Savings Type"
$elements.Range("L2").Value = "Customer-specific code for the type of third party savings"
